$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Student counts")

$ws.Range("B2").Value = 228
$ws.Range("B4").Value = 30
$ws.Range("B6").Value = 8
$ws.Range("B8").Value = 41
$ws.Range("B16").Value = 7
$ws.Range("B19").Value = 36
$ws.Range("B21").Value = 6
